# Helper: force a numeric-looking string to be stored as TEXT (shared string),
# matching the source workbook convention where every "value" column in
# Restricciones_del_follower / Punto_modificado / Vector_bf / Vector_BF is
# kept as text even though it looks like a plain number. Excel's COM layer
# auto-converts a bare numeric string typed into a General cell into a real
# number, so we briefly mark the cell as Text, assign, then clear the
# (now pointless) explicit formatting back off the cell.
function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower (sheet index 3) ---
$ws = $wb.Worksheets.Item(3)

$ws.Range("A2").Value = "-2.9128815680460503 + 0.15422855349071463y_1 + 0.6354648928176064y_2"
Set-TextValue $ws.Range("B2") "2.9128815680460503"
Set-TextValue $ws.Range("D2") "0.03"
Set-TextValue $ws.Range("E2") "3.0"
Set-TextValue $ws.Range("F2") "1.5"

$ws.Range("A3").Value = "-21.393257458149904 + 2.579393529463101y_1 + 0.8695410773557548y_2"
Set-TextValue $ws.Range("B3") "17.393257458149904"
Set-TextValue $ws.Range("D3") "0.85"
Set-TextValue $ws.Range("E3") "9.9"
Set-TextValue $ws.Range("F3") "9.1"

$ws.Range("A4").Value = "70.47052973353095 - 2x - 7.347310071313453y_1 - 0.5956431104918103y_2"
Set-TextValue $ws.Range("B4") "-86.47052973353095"
Set-TextValue $ws.Range("D4") "0.73"
Set-TextValue $ws.Range("E4") "8.8"
Set-TextValue $ws.Range("F4") "1.6"

$ws.Range("A5").Value = "-65.69063996964684 + 8x + 0.8257491743606533y_1 - 0.0959344506799863y_2"
Set-TextValue $ws.Range("B5") "17.000639969646844"
Set-TextValue $ws.Range("D5") "0.02"
Set-TextValue $ws.Range("E5") "2.0"
Set-TextValue $ws.Range("F5") "2.3000000000000003"

$ws.Range("A6").Value = "6.1925984835838275 - 2x + 0.6179237578818615y_1 + 1.4413078199944613y_2"
Set-TextValue $ws.Range("B6") "5.777401516416172"
Set-TextValue $ws.Range("D6") "0.72"
Set-TextValue $ws.Range("E6") "7.7"
Set-TextValue $ws.Range("F6") "4.2"

# --- Punto_modificado (sheet index 4) ---
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws.Range("A2") "7.4"
Set-TextValue $ws.Range("B2") "7.35"
Set-TextValue $ws.Range("C2") "2.8"

# --- Vector_bf (sheet index 5; name clashes case-insensitively with Vector_BF) ---
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("A2") "3.705004906248311"
Set-TextValue $ws.Range("A3") "-1.3591773332603108"

# --- Vector_BF (sheet index 6; name clashes case-insensitively with Vector_bf) ---
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "18.0"
Set-TextValue $ws.Range("A3") "35.248135740989916"
Set-TextValue $ws.Range("A4") "-18.17939328454424"

# --- Vector_Alpha (sheet index 7) -- these stay real numbers, not text ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 2.975681822770921
$ws.Range("A3").Value = 1.6382728748545337
